$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, even if it looks like a number,
# matching the source data (inline strings such as "1.00" or "3.422.54").
function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue 2 4 '61.999.39'
Set-TextValue 3 4 '3.422.54'
Set-TextValue 3 5 '  +1.02%  '
Set-TextValue 4 4 '1.00'
Set-TextValue 4 5 '  +0.02%  '
Set-TextValue 5 4 '578.96'
Set-TextValue 5 5 '  +1.22%  '
Set-TextValue 6 4 '144.81'
Set-TextValue 6 5 '  +2.06%  '
Set-TextValue 7 5 '  +0.07%  '
Set-TextValue 8 5 '  +0.09%  '
Set-TextValue 9 4 '7.63'
Set-TextValue 9 5 '  -0.37%  '
Set-TextValue 10 4 '0.124'
Set-TextValue 10 5 '  +0.68%  '
Set-TextValue 11 5 '  -0.26%  '
Set-TextValue 12 4 '4.008.36'
Set-TextValue 12 5 '  +1.06%  '
Set-TextValue 13 5 '  -0.67%  '
Set-TextValue 14 4 '28.27'
Set-TextValue 14 5 '  +1.55%  '
Set-TextValue 15 4 '3.424.79'
Set-TextValue 15 5 '  +1.49%  '
Set-TextValue 16 5 '  -0.30%  '
Set-TextValue 17 4 '62.033.23'
Set-TextValue 17 5 '  +1.53%  '
Set-TextValue 18 5 '  +1.53%  '
Set-TextValue 19 4 '13.94'
Set-TextValue 19 5 '  +1.98%  '
Set-TextValue 20 5 '  +3.29%  '
Set-TextValue 21 4 '390.60'
Set-TextValue 21 5 '  +1.99%  '
Set-TextValue 22 4 '74.50'
Set-TextValue 22 5 '  -1.05%  '
Set-TextValue 23 5 '  +0.18%  '
Set-TextValue 24 5 '  +0.28%  '
Set-TextValue 25 4 '0.0000116'
Set-TextValue 25 5 '  -0.42%  '
Set-TextValue 26 4 '0.192'
Set-TextValue 26 5 '  +4.77%  '
Set-TextValue 27 4 '7.47'
Set-TextValue 27 5 '  +3.10%  '
Set-TextValue 28 4 '0.999'
Set-TextValue 28 5 '  -0.01%  '
Set-TextValue 29 4 '8.05'
Set-TextValue 29 5 '  +1.04%  '
Set-TextValue 30 4 '2.16'
Set-TextValue 30 5 '  +0.60%  '
Set-TextValue 31 4 '1.43'
Set-TextValue 31 5 '  +3.10%  '
Set-TextValue 32 5 '  +0.01%  '
Set-TextValue 33 5 '  +1.36%  '
Set-TextValue 34 5 '  +6.02%  '
Set-TextValue 35 4 '6.99'
Set-TextValue 35 5 '  +0.39%  '
Set-TextValue 36 4 '168.54'
Set-TextValue 36 5 '  +1.18%  '
Set-TextValue 37 4 '3.455.78'
Set-TextValue 37 5 '  +1.04%  '
Set-TextValue 38 5 '  +0.73%  '
Set-TextValue 39 4 '28.57'
Set-TextValue 39 5 '  +7.42%  '
Set-TextValue 40 4 '0.0756'
Set-TextValue 40 5 '  -1.42%  '
Set-TextValue 41 4 '0.787'
Set-TextValue 41 5 '  +1.02%  '
Set-TextValue 42 5 '  +1.91%  '
Set-TextValue 43 5 '  +0.91%  '
Set-TextValue 44 4 '1.17'
Set-TextValue 44 5 '  +4.31%  '
Set-TextValue 45 4 '2.534.70'
Set-TextValue 45 5 '  +3.49%  '
Set-TextValue 46 4 '22.93'
Set-TextValue 46 5 '  -0.18%  '
Set-TextValue 47 2 'FirstDigitalUSD'
Set-TextValue 47 3 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 47 4 '1.00'
Set-TextValue 47 5 '  +0.04%  '
Set-TextValue 48 2 'Cosmos'
Set-TextValue 48 3 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 48 4 '6.62'
Set-TextValue 48 5 '  -0.92%  '
Set-TextValue 49 5 '  +0.26%  '
Set-TextValue 50 4 '2.10'
Set-TextValue 50 5 '  -2.56%  '
Set-TextValue 51 5 '  +0.08%  '
